# Weekly CompStat update: new crime data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Helper functions to change a cell's value while re-pointing its style ----
# donor cells (row 22, untouched by this week's edit) carry the three body-row
# styles we need to swap between: s14 (text placeholder), s15 (#,##0 integer),
# s16 (#,##0.0 percent-change). Copy + PasteSpecial(xlPasteFormats) re-points a
# cell at the existing style index instead of minting a new one.
function Set-StyledNumber($ref, $value, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $value
}

function Set-StyledText($ref, $text, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

$TextStyleDonor = "C22"   # s="14" (General, text placeholder)
$IntStyleDonor  = "J22"   # s="15" (#,##0)
$PctStyleDonor  = "K22"   # s="16" (#,##0.0;"-"#,##0.0)

# ---- Cells whose style flips between text-placeholder and numeric this week ----
Set-StyledNumber "D14" 1 $IntStyleDonor
Set-StyledNumber "E14" -100 $PctStyleDonor
Set-StyledNumber "G14" 1 $IntStyleDonor
Set-StyledNumber "H14" -100 $PctStyleDonor
Set-StyledNumber "J14" 1 $IntStyleDonor
Set-StyledNumber "K14" 0 $PctStyleDonor
Set-StyledText "D15" "0" $TextStyleDonor
Set-StyledText "E15" "***.*" $TextStyleDonor
Set-StyledText "D23" "0" $TextStyleDonor
Set-StyledText "E23" "***.*" $TextStyleDonor
Set-StyledText "D27" "0" $TextStyleDonor
Set-StyledText "E27" "***.*" $TextStyleDonor
Set-StyledNumber "C28" 4 $IntStyleDonor
Set-StyledText "D28" "0" $TextStyleDonor
Set-StyledText "E28" "***.*" $TextStyleDonor

# ---- Simple numeric value updates (style unchanged) ----
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = 65.217391304347
$ws.Range("L16").Value = 5.555555555555
$ws.Range("M16").Value = 65.217391304347
$ws.Range("N16").Value = -72.058823529411
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 31
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = -11.428571428571
$ws.Range("L17").Value = -18.421052631578
$ws.Range("M17").Value = 14.814814814814
$ws.Range("N17").Value = -35.416666666666
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = -18.75
$ws.Range("L18").Value = -7.142857142857
$ws.Range("M18").Value = -31.578947368421
$ws.Range("N18").Value = -88.793103448275
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 44.444444444444
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 28.571428571428
$ws.Range("I19").Value = 160
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = 32.231404958677
$ws.Range("L19").Value = 34.453781512605
$ws.Range("M19").Value = 162.295081967213
$ws.Range("N19").Value = 122.222222222222
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 93.75
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 24.285714285714
$ws.Range("L20").Value = 50
$ws.Range("M20").Value = 383.333333333333
$ws.Range("N20").Value = -73.952095808383
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 34.782608695652
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 31.111111111111
$ws.Range("I21").Value = 344
$ws.Range("J21").Value = 284
$ws.Range("K21").Value = 21.12676056338
$ws.Range("L21").Value = 19.8606271777
$ws.Range("M21").Value = 101.169590643275
$ws.Range("N21").Value = -58.604091456077
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 6.25
$ws.Range("L23").Value = -15
$ws.Range("M23").Value = 70
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -2.631578947368
$ws.Range("I24").Value = 215
$ws.Range("J24").Value = 204
$ws.Range("K24").Value = 5.392156862745
$ws.Range("L24").Value = 2.380952380952
$ws.Range("M24").Value = 38.709677419354
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 79
$ws.Range("J25").Value = 83
$ws.Range("K25").Value = -4.819277108433
$ws.Range("L25").Value = -34.710743801652
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 56.25
$ws.Range("I26").Value = 66
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = 17.857142857142
$ws.Range("L26").Value = -2.941176470588
$ws.Range("M26").Value = -5.714285714285
$ws.Range("G27").Value = 2
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = 37.5
$ws.Range("L28").Value = 57.142857142857
$ws.Range("N29").Value = -50
$ws.Range("N30").Value = -50

# ---- Masthead text: volume/issue number and report date range ----
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"
